$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 12:52"

# Target data for rows 4-65 (Ciudad, Casos totales, Casos activos, Recuperados, Muertes)
# The underlying per-province/city statistics refreshed to a newer snapshot, which
# re-sorts the leaderboard by "Casos totales" (column B) descending - so labels and
# values both shift between rows relative to the previous snapshot.

$rows = @(
    ,("Madrid", 42450, 19836, 17028, 5586)
    ,("Cataluña", 29647, 12250, 14356, 3041)
    ,("Galicia", 6331, 688, 5385, 258)
    ,("Bizkaia/Vizcaya", 4968, 4151, 4663, 317)
    ,("Ciudad Real", 4449, 1557, 8976, 400)
    ,("Valencia/Valencia", 3949, 733, 2873, 343)
    ,("Navarra", 3467, 450, 2811, 206)
    ,("Albacete", 3087, 1557, 8976, 263)
    ,("La Rioja", 2951, 1061, 1713, 177)
    ,("Araba/Alava", 2806, 4151, 4663, 229)
    ,("Alacant/Alicante", 2803, 469, 2042, 292)
    ,("Castilla-La Mancha", 2780, 71, 2446, 263)
    ,("Zaragoza", 2679, 484, 1932, 263)
    ,("Toledo", 2597, 1557, 8976, 352)
    ,("A Coruña", 1969, 333, 1788, 67)
    ,("Malaga", 1932, 289, 1518, 125)
    ,("Salamanca", 1807, 439, 1151, 217)
    ,("Sevilla", 1713, 105, 1496, 112)
    ,("Asturias", 1705, 254, 1349, 102)
    ,("Gipuzkoa/Guipuzcoa", 1678, 4151, 4663, 92)
    ,("Valladolid", 1602, 553, 896, 153)
    ,("Cantabria", 1572, 175, 1305, 92)
    ,("Granada", 1550, 182, 1240, 128)
    ,("Pontevedra", 1536, 333, 1411, 30)
    ,("Caceres", 1408, 113, 1075, 220)
    ,("Segovia", 1349, 387, 847, 115)
    ,("Leon", 1344, 564, 584, 196)
    ,("Murcia", 1326, 193, 1048, 85)
    ,("Tenerife", 1122, 249, 1422, 59)
    ,("Cordoba", 1055, 84, 932, 39)
    ,("Burgos", 1024, 387, 517, 120)
    ,("Guadalajara", 973, 1557, 8976, 133)
    ,("Jaen", 973, 62, 858, 53)
    ,("Aragon", 907, 29, 838, 40)
    ,("Castello/Castellon", 899, 142, 668, 89)
    ,("Cadiz", 881, 109, 740, 32)
    ,("Soria", 837, 172, 601, 64)
    ,("Avila", 766, 253, 430, 83)
    ,("Ourense", 751, 333, 660, 22)
    ,("Badajoz", 708, 169, 501, 38)
    ,("Cuenca", 682, 1557, 8976, 107)
    ,("Lugo", 586, 333, 520, 11)
    ,("Palencia", 510, 126, 343, 41)
    ,("Gran Canaria", 444, 249, 1422, 25)
    ,("Huesca", 428, 68, 309, 51)
    ,("Teruel", 406, 99, 273, 34)
    ,("Almeria", 371, 50, 298, 23)
    ,("Zamora", 342, 107, 196, 39)
    ,("Huelva", 292, 19, 257, 16)
    ,("Mallorca", 210, 18, 194, 12)
    ,("Melilla", 93, 12, 79, 2)
    ,("Ceuta", 84, 7, 73, 4)
    ,("La Palma", 67, 249, 1422, 3)
    ,("Lanzarote", 61, 249, 1422, 2)
    ,("Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena", 58, 0, 58, 3)
    ,("Fuerteventura", 23, 249, 1422, 0)
    ,("Ibiza", 21, 18, 20, 1)
    ,("Menorca", 15, 18, 13, 0)
    ,("La Gomera", 7, 249, 1422, 0)
    ,("Arroyo de la Luz", 7, 0, 7, 0)
    ,("El Hierro", 1, 249, 1422, 0)
    ,("Formentera", 0, 10, 0, 8)
)

$n = $rows.Count
$cols = 5
$arr = New-Object 'object[,]' $n,$cols
for ($i = 0; $i -lt $n; $i++) {
    $row = $rows[$i]
    for ($j = 0; $j -lt $cols; $j++) {
        $arr[$i,$j] = $row[$j]
    }
}

$ws.Range("A4:E65").Value = $arr
